# Update Cfp-Ncr1.xlsx LR-pairs sheet with refreshed TPM-derived numbers.
# The source data now only has 4 sending-cluster rows (rows 2-5); the old
# rows 6-9 (MuSCs / Resolving-Mac duplicate pairing) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-obsolete rows 6-9 first so the sheet shrinks to A1:T5.
$ws.Rows("6:9").Delete()

# Row 2 : ECs -> Cfp/Ncr1 -> Resolving-Mac
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cfp"
$ws.Range("C2").Value = "Ncr1"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.725679
$ws.Range("H2").Value = 5.177037
$ws.Range("I2").Value = 0.06358060230479941
$ws.Range("J2").Value = 0.06358060230479941
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2942353333333333
$ws.Range("N2").Value = 0.882706
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.5077557357913334
$ws.Range("R2").Value = 4.569801622122
$ws.Range("S2").Value = 0.06358060230479941
$ws.Range("T2").Value = 0.06358060230479941

# Row 3 : FAPs -> Cfp/Ncr1 -> Resolving-Mac
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cfp"
$ws.Range("C3").Value = "Ncr1"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.436639333333333
$ws.Range("H3").Value = 7.309918
$ws.Range("I3").Value = 0.0897750951439394
$ws.Range("J3").Value = 0.0897750951439394
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2942353333333333
$ws.Range("N3").Value = 0.882706
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.7169453864564445
$ws.Range("R3").Value = 6.452508478107999
$ws.Range("S3").Value = 0.0897750951439394
$ws.Range("T3").Value = 0.0897750951439394

# Row 4 : MuSCs -> Cfp/Ncr1 -> Resolving-Mac
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Cfp"
$ws.Range("C4").Value = "Ncr1"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.158792
$ws.Range("H4").Value = 6.476376
$ws.Range("I4").Value = 0.07953813867514323
$ws.Range("J4").Value = 0.07953813867514324
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2942353333333333
$ws.Range("N4").Value = 0.882706
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.6351928837173334
$ws.Range("R4").Value = 5.716735953456
$ws.Range("S4").Value = 0.07953813867514323
$ws.Range("T4").Value = 0.07953813867514324

# Row 5 : Resolving-Mac -> Cfp/Ncr1 -> Resolving-Mac
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Cfp"
$ws.Range("C5").Value = "Ncr1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 20.82048533333333
$ws.Range("H5").Value = 62.461456
$ws.Range("I5").Value = 0.767106163876118
$ws.Range("J5").Value = 0.767106163876118
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2942353333333333
$ws.Range("N5").Value = 0.882706
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 6.126122442215111
$ws.Range("R5").Value = 55.135101979936
$ws.Range("S5").Value = 0.767106163876118
$ws.Range("T5").Value = 0.767106163876118
